# Add "undo" option handling: swap the "Taken By" / "Received" headers
# and adjust several rows' Taken-By / Received values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: swap C1/D1 ---
$ws.Range("C1").Value = "Received"
$ws.Range("D1").Value = "Taken By"

# --- Row 2 ---
$ws.Range("C2").Value = "Yes"
$ws.Range("D2").Value = "BOB"

# --- Row 3 ---
$ws.Range("D3").Value = ""

# --- Row 4 ---
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""

# --- Row 5 ---
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""

# --- Row 7 ---
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""

# --- Row 8 ---
$ws.Range("C8").Value = "o"

# --- Row 11 ---
$ws.Range("C11").Value = "o"

# --- Row 30 ---
$ws.Range("D30").Value = ""

# --- Row 31 ---
$ws.Range("C31").Value = ""
$ws.Range("D31").Value = ""

# --- Row 90 ---
$ws.Range("C90").Value = "o"
